$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (D, M, N, O, P, S columns) got reshuffled between rows as
# part of a weekly re-sort of the "Femacal de La Calera - Breva" records.
# Capture the "before" values for each affected row first, then write the
# "after" values, so that rows reading from one another don't clobber
# data before it has been captured.

$rows = @(2, 3, 4, 5, 7, 9)

$before = @{}
foreach ($r in $rows) {
    $before[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        S = $ws.Cells.Item($r, 19).Value2
    }
}

# Row r takes the previously-captured values of row srcMap[r]
$srcMap = @{
    2 = 5
    3 = 4
    4 = 7
    5 = 3
    7 = 9
    9 = 2
}

foreach ($r in $rows) {
    $src = $before[$srcMap[$r]]
    $ws.Cells.Item($r, 4).Value2 = $src.D
    $ws.Cells.Item($r, 13).Value2 = $src.M
    $ws.Cells.Item($r, 14).Value2 = $src.N
    $ws.Cells.Item($r, 15).Value2 = $src.O
    $ws.Cells.Item($r, 16).Value2 = $src.P
    $ws.Cells.Item($r, 19).Value2 = $src.S
}
